$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new values in column B, rows 3 and 4
# (written in this order so the shared-string table indices match the target: 5="8.81100010872 Seconds", 6="194.180000000 seconds")
$ws.Range("B4").Value = "8.81100010872 Seconds"
$ws.Range("B3").Value = "194.180000000 seconds"

# Update the active selection to C1
$ws.Range("C1").Select()
